# Update "想去人数" (want-to-go count, column F) for several events across
# the "展览", "演出" and "全部类型" sheets, matching the site's regenerated
# output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1533
$ws.Range("F3").Value  = 879
$ws.Range("F6").Value  = 514
$ws.Range("F7").Value  = 7716
$ws.Range("F11").Value = 5572
$ws.Range("F14").Value = 7693
$ws.Range("F15").Value = 9095
$ws.Range("F17").Value = 911
$ws.Range("F18").Value = 4479
$ws.Range("F26").Value = 1679
$ws.Range("F28").Value = 944
$ws.Range("F32").Value = 2316
$ws.Range("F35").Value = 1474
$ws.Range("F40").Value = 2985
$ws.Range("F41").Value = 4119
$ws.Range("F44").Value = 427
$ws.Range("F48").Value = 178
$ws.Range("F49").Value = 4093

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F26").Value = 99

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 1533
$ws.Range("F4").Value  = 879
$ws.Range("F7").Value  = 514
$ws.Range("F11").Value = 5572
$ws.Range("F13").Value = 7693
$ws.Range("F16").Value = 911
$ws.Range("F25").Value = 1679
$ws.Range("F27").Value = 944
$ws.Range("F31").Value = 2316
$ws.Range("F40").Value = 4119
$ws.Range("F44").Value = 427
$ws.Range("F48").Value = 178
$ws.Range("F49").Value = 4093
